$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Make the first paragraph ("8.3.1") bold.
# ------------------------------------------------------------------
$d.Paragraphs.Item(1).Range.Font.Bold = $true

# ------------------------------------------------------------------
# 2) Move the (single, Word-managed) "_GoBack" bookmark so it now
#    sits right after the run "Come in! Hi Khae" (end of that
#    paragraph, before the pilcrow) instead of its old spot after
#    "Jim can talk to Anne in English". Re-adding a bookmark with an
#    existing name relocates it, which also takes care of removing
#    it from its old location.
#    The engine snaps a bookmark collapsed exactly at
#    "paragraph-end-1" back to the start of the paragraph, so we
#    temporarily insert a marker character after the text, add the
#    bookmark right before that marker (which is no longer a
#    paragraph boundary), then remove the marker again.
# ------------------------------------------------------------------
$comeRange = $d.Content
$comeRange.Find.Execute("Come in! Hi Khae", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$comeEnd = $comeRange.End
$comeRange.InsertAfter("X")
$bmRange = $d.Range($comeEnd, $comeEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
$markerRange = $d.Range($comeEnd, $comeEnd + 1)
$markerRange.Text = ""

# ------------------------------------------------------------------
# 3) Merge the two runs "No, you can't. " and
#    "The Big Boss says you can't watch TV at work!" into one run.
# ------------------------------------------------------------------
$mergedText = "No, you can" + [char]0x2019 + "t. The Big Boss says you can" + [char]0x2019 + "t watch TV at work!"
$mergeRange = $d.Content
$mergeRange.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeRange.Text = "TEMP_PLACEHOLDER_8F3"
$mergeRange.Text = $mergedText
